$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    # Several "Price" cells hold numeric-looking text (e.g. "316.30", "1.003") that must
    # stay exactly as authored (fixed decimal places, thousands separated by dots, etc.).
    # A plain Range.Value assignment lets Excel auto-convert such strings into real
    # numbers (dropping trailing zeros) and, if forced via NumberFormat/quote-prefix,
    # stamps the cell with a new style index that was not present in the source file.
    # Writing the text as a formula result and collapsing it to a static value with
    # Copy + PasteSpecial(xlPasteValues) keeps both the exact text and the original
    # (unstyled) cell formatting intact.
    $r = $ws.Range($cellRef)
    $r.Formula = '="' + $value + '"'
    $r.Copy()
    $r.PasteSpecial(-4163)
}

$ws.Range("D2").Value = '24.792.82'
$ws.Range("E2").Value = '  +0.61%  '
$ws.Range("D3").Value = '1.702.68'
$ws.Range("E3").Value = '  +0.12%  '
$ws.Range("E4").Value = '  +0.33%  '
Set-TextValue "D5" "316.30"
$ws.Range("E5").Value = '  +0.12%  '
Set-TextValue "D6" "1.003"
$ws.Range("E6").Value = '  +0.32%  '
Set-TextValue "D7" "0.3936"
$ws.Range("E7").Value = '  -0.21%  '
Set-TextValue "D8" "0.4042"
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("E9").Value = '  -2.37%  '
Set-TextValue "D10" "53.94"
$ws.Range("E10").Value = '  -1.69%  '
Set-TextValue "D11" "1.004"
$ws.Range("E11").Value = '  +0.37%  '
Set-TextValue "D12" "0.08907"
$ws.Range("E12").Value = '  +0.95%  '
Set-TextValue "D13" "7.248"
$ws.Range("E13").Value = '  -0.78%  '
Set-TextValue "D14" "23.44"
$ws.Range("E14").Value = '  -0.01%  '
Set-TextValue "D15" "8.014"
$ws.Range("E15").Value = '  +4.82%  '
Set-TextValue "D16" "0.00001330"
$ws.Range("E16").Value = '  -0.40%  '
$ws.Range("D17").Value = '1.701.98'
$ws.Range("E17").Value = '  -0.48%  '
Set-TextValue "D18" "100.15"
$ws.Range("E18").Value = '  -0.67%  '
Set-TextValue "D19" "0.07043"
$ws.Range("E19").Value = '  -0.34%  '
Set-TextValue "D20" "19.67"
$ws.Range("E20").Value = '  -0.28%  '
Set-TextValue "D21" "7.034"
$ws.Range("E21").Value = '  +1.28%  '
Set-TextValue "D22" "1.002"
$ws.Range("E22").Value = '  +0.14%  '
Set-TextValue "D23" "14.59"
$ws.Range("E23").Value = '  +2.89%  '
$ws.Range("D24").Value = '24.784.87'
$ws.Range("E24").Value = '  +0.64%  '
Set-TextValue "D25" "3.203"
$ws.Range("E25").Value = '  +7.22%  '
Set-TextValue "D26" "2.355"
$ws.Range("E26").Value = '  +1.16%  '
$ws.Range("E27").Value = '  +1.61%  '
Set-TextValue "D28" "162.06"
$ws.Range("E28").Value = '  +1.15%  '
Set-TextValue "D29" "136.49"
$ws.Range("E29").Value = '  +1.84%  '
$ws.Range("B30").Value = 'Filecoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue "D30" "7.943"
$ws.Range("E30").Value = '  +2.00%  '
$ws.Range("B31").Value = 'HuobiToken'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue "D31" "5.173"
$ws.Range("E31").Value = '  -1.12%  '
Set-TextValue "D32" "0.08781"
$ws.Range("E32").Value = '  +2.32%  '
Set-TextValue "D33" "1.082"
$ws.Range("E33").Value = '  -2.92%  '
Set-TextValue "D34" "7.236"
$ws.Range("E34").Value = '  -2.95%  '
Set-TextValue "D35" "11.22"
$ws.Range("E35").Value = '  +0.22%  '
Set-TextValue "D36" "1.977"
$ws.Range("E36").Value = '  +0.96%  '
Set-TextValue "D37" "0.2741"
$ws.Range("E37").Value = '  -0.90%  '
Set-TextValue "D38" "14.40"
$ws.Range("E38").Value = '  -2.82%  '
Set-TextValue "D39" "0.09200"
$ws.Range("E39").Value = '  +1.60%  '
Set-TextValue "D40" "0.02752"
$ws.Range("E40").Value = '  -1.20%  '
Set-TextValue "D41" "1.461"
$ws.Range("E41").Value = '  -0.73%  '
Set-TextValue "D42" "0.7705"
$ws.Range("E42").Value = '  -1.06%  '
Set-TextValue "D43" "15.81"
$ws.Range("E43").Value = '  +0.27%  '
Set-TextValue "D44" "0.7172"
$ws.Range("E44").Value = '  -1.79%  '
$ws.Range("E45").Value = '  +2.21%  '
Set-TextValue "D46" "4.210"
$ws.Range("E46").Value = '  +0.04%  '
Set-TextValue "D47" "1.003"
$ws.Range("E47").Value = '  +0.26%  '
Set-TextValue "D48" "140.65"
$ws.Range("E48").Value = '  -1.04%  '
Set-TextValue "D49" "1.316"
$ws.Range("E49").Value = '  -1.60%  '
Set-TextValue "D50" "90.95"
$ws.Range("E50").Value = '  +2.76%  '
Set-TextValue "D51" "0.07996"
$ws.Range("E51").Value = '  -0.59%  '
